$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Map of cell address -> new value, derived from the XML diff.
$updates = @{
    "O6" = 1.5
    "P6" = 2.5
    "U13" = 1.47
    "Q14" = 2.4
    "R14" = 1.53
    "M15" = 1.06
    "N15" = 10
    "Q15" = 2.03
    "R15" = 1.78
    "M17" = 1.06
    "N17" = 10
    "M22" = 1.07
    "O22" = 1.36
    "M23" = 1.07
    "O23" = 1.36
    "M24" = 1.06
    "O24" = 1.33
    "M25" = 1.08
    "O25" = 1.36
    "K27" = 2.07
    "S27" = 1.42
    "T27" = 2.67
    "U27" = 1.75
    "V27" = 1.95
    "X27" = 11
    "Y27" = 8.75
    "AB27" = 28
    "AT27" = 2.67
    "AZ27" = 90
    "G28" = 2.4
    "H28" = 3.15
    "K28" = 2.05
    "L28" = 3.4
    "W28" = 8.5
    "X28" = 12.5
    "AA28" = 19.5
    "AB28" = 27
    "AG28" = 9
    "AL28" = 30
    "AM28" = 400
    "AX28" = 15.5
    "AY28" = 23
    "BA28" = 110
    "BB28" = 300
    "G34" = 2.6
    "H34" = 3.2
    "I34" = 2.47
    "K34" = 2.12
    "L34" = 3.1
    "M34" = 1.07
    "N34" = 7
    "O34" = 1.33
    "P34" = 3.05
    "S34" = 1.4
    "U34" = 1.75
    "V34" = 1.95
    "W34" = 8.5
    "X34" = 13
    "Y34" = 9.75
    "Z34" = 29
    "AC34" = 7
    "AD34" = 6.3
    "AE34" = 14
    "AF34" = 65
    "AJ34" = 27
    "AM34" = 500
    "AN34" = 4.6
    "AO34" = 14
    "AP34" = 21
    "AR34" = 90
    "AS34" = 250
    "AU34" = 7
    "AV34" = 60
    "AW34" = 4.45
    "AX34" = 13.5
    "AZ34" = 60
    "P35" = 2.87
    "Q35" = 1.91
    "R35" = 1.7
    "V35" = 1.8
    "W35" = 6.5
    "AL35" = 50
    "AP35" = 17.5
    "AV35" = 65
    "M36" = 1.01
    "N36" = 7.8
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

Write-Output ("Updated {0} cells." -f $updates.Count)
